$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Row 2: Giovani / 0803 / Galpão Toyota / ... / Pendente / Maxvel: 36 / Forte: 13
$ws.Range("A2").Value = "Giovani"
$ws.Range("B2").Value = "'0803"
$ws.Range("C2").Value = "Galpão Toyota"
$ws.Range("D2").Value = "Pegar MAC da central pra base,  instalar switch, passar central pra internet e passar acesso das câmeras para a base."
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = "Pendente"
$ws.Range("H2").Value = "Maxvel: 36 / Forte: 13"

# Row 3: Giovani / 0701 / Usina Amaral / Restaurar comunicação geral da usina. / Pendente
$ws.Range("A3").Value = "Giovani"
$ws.Range("B3").Value = "'0701"
$ws.Range("C3").Value = "Usina Amaral"
$ws.Range("D3").Value = "Restaurar comunicação geral da usina."
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = ""
$ws.Range("G3").Value = "Pendente"

# Row 4: Giovani / 0865 / MW Educação / Pegar MAC da central para a base poder fazer uma alteração. / Pendente
$ws.Range("A4").Value = "Giovani"
$ws.Range("B4").Value = "'0865"
$ws.Range("C4").Value = "MW Educação"
$ws.Range("D4").Value = "Pegar MAC da central para a base poder fazer uma alteração."
$ws.Range("G4").Value = "Pendente"

# Row 5: Giovani / 0372 / Cmei José Clementino / Várias câmeras fora, parece ser problema de fonte. / Pendente
$ws.Range("A5").Value = "Giovani"
$ws.Range("B5").Value = "'0372"
$ws.Range("C5").Value = "Cmei José Clementino"
$ws.Range("D5").Value = "Várias câmeras fora, parece ser problema de fonte."
$ws.Range("G5").Value = "Pendente"

# Update view: top-left cell and selection
$ws.Application.ActiveWindow.ScrollColumn = 5
$ws.Range("G5").Select()
